$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'49.803.31"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "'2.642.55"
$ws.Range("E3").Value = "  +5.68%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'113.39"
$ws.Range("E5").Value = "  +6.49%  "
$ws.Range("D6").Value = "'326.50"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.552"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "'40.93"
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "'0.0818"
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'7.32"
$ws.Range("E14").Value = "  +3.41%  "
$ws.Range("D15").Value = "'3.060.24"
$ws.Range("E15").Value = "  +5.85%  "
$ws.Range("D16").Value = "'2.646.71"
$ws.Range("E16").Value = "  +6.00%  "
$ws.Range("D17").Value = "'0.868"
$ws.Range("E17").Value = "  +4.07%  "
$ws.Range("D18").Value = "'49.699.76"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").Value = "'13.10"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.73"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "'2.92"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'0.0₃0952"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "'72.07"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "'276.85"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").Value = "'26.68"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'9.98"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").Value = "'35.88"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'50.35"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("D33").Value = "'19.50"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "'0.0803"
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'2.06"
$ws.Range("E37").Value = "  +6.18%  "
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("E39").Value = "  +6.13%  "
$ws.Range("D40").Value = "'124.62"
$ws.Range("E40").Value = "  +2.75%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'22.03"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'0.0313"
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'2.066.53"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.32"
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("D47").Value = "'2.29"
$ws.Range("E47").Value = "  +14.19%  "
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("D49").Value = "'9.07"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").Value = "'5.36"
$ws.Range("E50").Value = "  +3.65%  "
$ws.Range("D51").Value = "'58.79"
$ws.Range("E51").Value = "  +3.56%  "
